# "on First Vote Notification" -- add a new To Do item that records the
# follow-up task for the notifications table migration, and add a "Done"
# status column to the Notifications sheet for the new row (On First Vote).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "To Do": append a new task row (row 30)
# ---------------------------------------------------------------------
$todo = $wb.Worksheets.Item(1)

$todo.Cells.Item(30, 1).Value = "migrate data_id in notifications table"
$todo.Cells.Item(30, 2).Value = "Done"

# Move the frozen-pane selection down one row, same as the authored sheet.
$todo.Activate() | Out-Null
$todo.Range("A31").Select() | Out-Null

# ---------------------------------------------------------------------
# Sheet "Notifications": add a "Done" status column (E) and mark the
# "On First Vote" notification row (row 2) as Done.
# ---------------------------------------------------------------------
$notif = $wb.Worksheets.Item(2)

$e2 = $notif.Range("E2")
$e2.Value = "Done"
$e2.VerticalAlignment = -4108     # xlCenter
$e2.HorizontalAlignment = -4108   # xlCenter

# Give the new column a width close to its authored size.
$notif.Columns.Item(5).ColumnWidth = 15.6

$notif.Activate() | Out-Null
$notif.Range("C3").Select() | Out-Null
